# Technical Description of System - apply edits per commit diff.
$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper pattern used repeatedly below: to insert a new bold heading
# paragraph right before an existing paragraph located by Find, we
# - find the anchor text,
# - remember the (pre-insert) paragraph index of that anchor - this is
#   exactly the index the freshly inserted blank paragraph will occupy,
# - call InsertParagraphBefore() on a collapsed range at the start of the
#   anchor paragraph,
# - fill that now-existing blank paragraph with the heading text and turn
#   on bold for its (whole) run.
# -----------------------------------------------------------------------

# =======================================================================
# 1) "Actions" heading before "The system begins by asking the user..."
# =======================================================================
$anchor = $d.Content
$anchor.Find.Execute("The system begins by asking the user") | Out-Null
$headingIndex = $anchor.Paragraphs(1).Index
$headStart = $anchor.Duplicate
$headStart.Collapse(1)
$headStart.InsertParagraphBefore() | Out-Null
$d.Paragraphs($headingIndex).Range.Text = "Actions"
$d.Paragraphs($headingIndex).Range.Font.Bold = 1

# =======================================================================
# 2) After the "Actions" paragraph (now the last paragraph in the body),
#    append: blank paragraph, "Design" heading, design paragraph.
# =======================================================================
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter() | Out-Null

$blankPara = $d.Paragraphs($d.Paragraphs.Count)
$blankPara.Range.InsertParagraphAfter() | Out-Null

$designHeadingIndex = $d.Paragraphs.Count
$d.Paragraphs($designHeadingIndex).Range.Text = "Design"
$d.Paragraphs($designHeadingIndex).Range.Font.Bold = 1

$designHeadingPara = $d.Paragraphs($d.Paragraphs.Count)
$designHeadingPara.Range.InsertParagraphAfter() | Out-Null

$designPara = $d.Paragraphs($d.Paragraphs.Count)
$designPara.Range.InsertAfter("The program is simple and easy to use, meaning the user does not need extra knowledge to navigate it. All options are clearly labeled and intuitive as well reversible. Information given is complete to the point where it doesn’t leave the user wondering what is trying to be explained, but enough white space is present so that it doesn’t appear cluttered. These design choices will ensure that the max number of people can use the program to its full potential. ") | Out-Null

# Split the design paragraph's single run into three runs (matching the
# source, which has each sentence as its own <w:r>) by bookmarking and
# immediately un-bookmarking the two split points - Word always splits
# the run at a bookmark edge, and deleting the bookmark again leaves the
# split in place without any bookmark residue.
$split1 = $d.Content
$split1.Find.Execute("reversible. ") | Out-Null
$split1.Collapse(0)
$split1.Bookmarks.Add("zzSplit1") | Out-Null
$d.Bookmarks("zzSplit1").Delete()

$split2 = $d.Content
$split2.Find.Execute("cluttered. ") | Out-Null
$split2.Collapse(0)
$split2.Bookmarks.Add("zzSplit2") | Out-Null
$d.Bookmarks("zzSplit2").Delete()

# =======================================================================
# 3) "Languages/Services" heading before "The frontend portion..."
# =======================================================================
$anchor = $d.Content
$anchor.Find.Execute("The frontend portion of this system") | Out-Null
$headingIndex = $anchor.Paragraphs(1).Index
$headStart = $anchor.Duplicate
$headStart.Collapse(1)
$headStart.InsertParagraphBefore() | Out-Null
$d.Paragraphs($headingIndex).Range.Text = "Languages/Services"
$d.Paragraphs($headingIndex).Range.Font.Bold = 1

# Swap the trailing sentence of that paragraph for the new sentence about
# the Pizza object (this also removes the old _GoBack bookmark that used
# to sit right after it - Word only ever keeps a single _GoBack bookmark,
# and step 5 below adds a fresh one elsewhere).
$d.Content.Find.Execute("The system allows all information entered to be changed or updated and all choices to be undone. ", $false, $false, $false, $false, $false, $true, 1, $false, "An object called Pizza is used by the program to save data into the database. ", 2) | Out-Null

# =======================================================================
# 4) "Purpose" heading before "The purpose of this system..."
# =======================================================================
$anchor = $d.Content
$anchor.Find.Execute("The purpose of this system") | Out-Null
$headingIndex = $anchor.Paragraphs(1).Index
$headStart = $anchor.Duplicate
$headStart.Collapse(1)
$headStart.InsertParagraphBefore() | Out-Null
$d.Paragraphs($headingIndex).Range.Text = "Purpose"
$d.Paragraphs($headingIndex).Range.Font.Bold = 1

# =======================================================================
# 5) Blank paragraph right under the title loses its center alignment.
# =======================================================================
$d.Paragraphs(2).Alignment = 0

# =======================================================================
# 6) Title: split "Technical Description of System" into two runs with
#    the _GoBack bookmark in between (Word keeps only a single _GoBack
#    bookmark, so adding it here removes the one from step 3 above).
# =======================================================================
$titleSplit = $d.Content
$titleSplit.Find.Execute("Technical Des") | Out-Null
$titleSplit.Collapse(0)
$titleSplit.Bookmarks.Add("_GoBack") | Out-Null
